$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New swap-request rows to append below the existing data (rows 14-21)
$newRows = @(
    @(13, "Vishal Dhanasekaran", 15, "e",                    "Mar 20", "G", "M", "Accepted", "2024-05-21 16:02:47"),
    @(13, "Vishal Dhanasekaran", 12, "b",                    "Mar 20", "M", "A", "Accepted", "2024-05-21 18:23:41"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "M", "N", "Accepted", "2024-05-22 17:21:01"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "N", "M", "Accepted", "2024-05-22 17:25:05"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "M", "N", "Declined", "2024-05-22 17:29:04"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "M", "N", "Accepted", "2024-05-22 18:23:08"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "N", "M", "Accepted", "2024-05-22 22:08:40"),
    @(13, "Vishal Dhanasekaran", 11, "Vishal Dhanasekaran",  "Mar 1",  "M", "N", "Accepted", "2024-05-22 22:27:07")
)

$startRow = 14
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
}
